# Append a new data row (row 54) to the active worksheet, mirroring the
# most recent temperature reading row (row 53), per data pulled from
# Adafruit IO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 54

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "25" looks numeric, so a plain assignment would be auto-coerced into a
# number by Excel. Force it to stay text (matching the source column,
# which is entirely text) using a leading apostrophe, then reset the
# cell style back to Normal so no stray quote-prefix formatting lingers.
$ws.Cells.Item($row, 3).Value = "'25"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
